$d = $word.ActiveDocument

# Split "Statement of Requirements " into "Statement of Requireme" + "nts "
# by placing the cursor between "Requireme" and "nts" (simulating a user
# clicking there and typing), which is what moves the _GoBack bookmark.

$range = $d.Content
$range.Find.Execute("Statement of Requirements ", $true, $false, $false, $false, $false, $true, 1, $false, "Statement of Requireme", 2)
